# Apply the recorded changes to the workbook.
$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("Input")
$wsRepay  = $wb.Worksheets.Item("Repayment Schedule")
$wsTrans  = $wb.Worksheets.Item("Transactions")

# -- Input sheet: it is no longer the selected/active tab -> deselect it by
#    activating a different sheet later (tabSelected gets recomputed on save).

# -- Repayment Schedule sheet: move the blank formatted cell from P2 to O2
#    (same style, empty contents), then scroll/re-point the view.
$wsRepay.Activate()
$wsRepay.Range("P2").Copy($wsRepay.Range("O2"))
$wsRepay.Range("P2").Clear()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$wsRepay.Range("C8").Select()

# -- Transactions sheet: update data values, column width, and drop column K.
$wsTrans.Range("A2").Value = 27
$wsTrans.Range("F2").Value = 785.8
$wsTrans.Range("G2").Value = 101.92
$wsTrans.Range("J2").Value = 9214.2
$wsTrans.Range("A3").Value = 26

$wsTrans.Range("K2").Clear()
$wsTrans.Columns("C").ColumnWidth = 15.65

# -- Transactions becomes the active (selected) tab / sheet.
$wsTrans.Activate()
$wsTrans.Range("C15").Select()

$wb.Save()
